# The data row for colo "XNN" (Xining, China) at worksheet row 268 was
# removed from the source dataset. Removing it shifts every following
# row up by one (e.g. the old row 269 "FRU" becomes the new row 268,
# and the former last row 333 "YWG" becomes the new last row 332), and
# the used range shrinks from A1:H333 to A1:H332.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 268 (colo "XNN" / Xining, China). This shifts all
# subsequent rows up by one, exactly matching the diff.
$ws.Rows.Item(268).Delete()
